$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sales by advisor/client per product category)
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 31 - LINDAO ZUÑIGA BRYAN JOSE / GUZMAN MAYORGA ROMINA SISNEY
$wsGrupo.Range("D31").Value = 2280.96   # 240X80 PORCELANATO
$wsGrupo.Range("L31").Value = 1961.56   # PIEDRA SINTERIZADA
$wsGrupo.Range("M31").Value = 114.61    # PORCELANATO

# Row 60 - totals "<n> de 58" per column (counts of non-zero rows)
$wsGrupo.Range("D60").Value = "2 de 58"
$wsGrupo.Range("L60").Value = "2 de 58"
$wsGrupo.Range("M60").Value = "5 de 58"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (sales by advisor/client per month)
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 31 - LINDAO ZUÑIGA BRYAN JOSE / GUZMAN MAYORGA ROMINA SISNEY
$wsMensual.Range("F31").Value = 4357.13   # noviembre

# Row 60 - totals
$wsMensual.Range("F60").Value = 33985.7

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (monthly compliance per product category)
# ---------------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 - 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 2738.88
$wsCumpl.Range("E3").Value = -871.1900000000001
$wsCumpl.Range("F3").Value = 1.466453212256852

# Row 11 - PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 8034.04
$wsCumpl.Range("E11").Value = -6593.12
$wsCumpl.Range("F11").Value = 5.575632234960997

# Row 12 - PORCELANATO
$wsCumpl.Range("D12").Value = 21123.07
$wsCumpl.Range("E12").Value = 26917.93
$wsCumpl.Range("F12").Value = 0.4396883911658792

# Row 14 - TOTAL
$wsCumpl.Range("D14").Value = 34097.58
$wsCumpl.Range("E14").Value = 23789.77196497848
$wsCumpl.Range("F14").Value = 0.5890333353066287
